$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 1).Value = "2025-10-28 18:36:33"
}
